$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first worksheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 76
$ws1.Range("F3").Value = 658
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F6").Value = 9961
$ws1.Range("F7").Value = 901
$ws1.Range("F10").Value = 4998
$ws1.Range("F11").Value = 4
$ws1.Range("F12").Value = 6
$ws1.Range("F15").Value = 63
$ws1.Range("F17").Value = 294
$ws1.Range("F18").Value = 575
$ws1.Range("F19").Value = 114
$ws1.Range("F21").Value = 7
$ws1.Range("F22").Value = 1501

# Sheet "全部类型" (All types) - fourth worksheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 76
$ws4.Range("F4").Value = 658
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F7").Value = 9961
$ws4.Range("F8").Value = 901
$ws4.Range("F11").Value = 4998
$ws4.Range("F12").Value = 4
$ws4.Range("F13").Value = 6
$ws4.Range("F16").Value = 63
$ws4.Range("F18").Value = 294
$ws4.Range("F19").Value = 575
$ws4.Range("F20").Value = 114
$ws4.Range("F22").Value = 7
$ws4.Range("F23").Value = 1501
